# Apply LIPID MAPS converter test-data changes to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E gets a Text ("@") number format (adds a new cellXfs entry). ---
$ws.Range("E1:E17").NumberFormat = "@"

# --- Column E: two new rows of data are inserted at E2:E3, pushing the
# --- previously-existing E2:E16 values down to E4:E18, plus one brand new
# --- value appended at E18. ---
$ws.Range("E2").Value  = "16:0"
$ws.Range("E3").Value  = "20:2(10Z,13E)(9Ke,15OH[S])"
$ws.Range("E4").Value  = "20:4(5Z,9E,11Z,14Z)(8OH)"
$ws.Range("E5").Value  = "22:5(8E,10Z,13Z,16Z,19Z)(7OH)"
$ws.Range("E6").Value  = "PE(16:0/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E7").Value  = "PE(16:0/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E8").Value  = "PC(16:0/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E9").Value  = "PC(16:0/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E10").Value = "PE(18:0/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E11").Value = "PE(18:0/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E12").Value = "PC(18:0/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E13").Value = "PC(18:0/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E14").Value = "PE(18:1/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E15").Value = "PE(18:1/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E16").Value = "PC(18:1/20:4(5Z,9E,11Z,14Z)(8OH))"
$ws.Range("E17").Value = "PC(18:1/22:5(8E,10Z,13Z,16Z,19Z)(7OH))"
$ws.Range("E18").Value = "BAD_Test5"

# --- Two "BAD_Test" values land in row 44 (columns C/D), extending that row. ---
$ws.Range("C44").Value = "BAD_Test3"
$ws.Range("D44").Value = "BAD_Test4"

# --- New trailing row 53 holds the remaining two "BAD_Test" values. ---
$ws.Range("A53").Value = "BAD_Test1"
$ws.Range("B53").Value = "BAD_Test2"

# --- Selection / active cell moves to E18 (matches the authored edit). ---
$ws.Range("E18").Select()
